$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 168 (Indice 167): Otelul 0 - 2 FC Botosani ---
$ws.Range("A168").Value = 167
$ws.Range("B168").Value = "romania"
$ws.Range("C168").Value = "liga-1"
$ws.Range("D168").Value = "2023-2024"
$ws.Range("E168").Value = 45282.6875
$ws.Range("F168").Value = "Otelul"
$ws.Range("G168").Value = 0
$ws.Range("H168").Value = "FC Botosani"
$ws.Range("I168").Value = 2
$ws.Range("J168").Value = 1.85
$ws.Range("K168").Value = "18/12/2023 19:12"
$ws.Range("L168").Value = 1.78
$ws.Range("M168").Value = "22/12/2023 16:29"
$ws.Range("N168").Value = 3.23
$ws.Range("O168").Value = "18/12/2023 19:12"
$ws.Range("P168").Value = 3.37
$ws.Range("Q168").Value = "22/12/2023 16:29"
$ws.Range("R168").Value = 4.37
$ws.Range("S168").Value = "18/12/2023 19:12"
$ws.Range("T168").Value = 5.3
$ws.Range("U168").Value = "22/12/2023 16:27"
$ws.Range("V168").Value = "https://www.betexplorer.com/football/romania/liga-1/otelul-fc-botosani/x0QRYveE/"

# --- Row 169 (Indice 168): Din. Bucuresti 1 - 0 FC Voluntari ---
$ws.Range("A169").Value = 168
$ws.Range("B169").Value = "romania"
$ws.Range("C169").Value = "liga-1"
$ws.Range("D169").Value = "2023-2024"
$ws.Range("E169").Value = 45282.8125
$ws.Range("F169").Value = "Din. Bucuresti"
$ws.Range("G169").Value = 1
$ws.Range("H169").Value = "FC Voluntari"
$ws.Range("I169").Value = 0
$ws.Range("J169").Value = 2.7
$ws.Range("K169").Value = "18/12/2023 19:12"
$ws.Range("L169").Value = 2.39
$ws.Range("M169").Value = "22/12/2023 19:01"
$ws.Range("N169").Value = 3.01
$ws.Range("O169").Value = "18/12/2023 19:12"
$ws.Range("P169").Value = 2.9
$ws.Range("Q169").Value = "22/12/2023 19:01"
$ws.Range("R169").Value = 2.7
$ws.Range("S169").Value = "18/12/2023 19:12"
$ws.Range("T169").Value = 3.56
$ws.Range("U169").Value = "22/12/2023 19:01"
$ws.Range("V169").Value = "https://www.betexplorer.com/football/romania/liga-1/din-bucuresti-voluntari/raleHCPt/"

# Copy formatting (styles only) from the last existing data row (167) so that
# the new rows 168-169 match the look of the rest of the table: bold/bordered
# centered index column (A) and a date-time number format on column E.
$ws.Range("A167:V167").Copy()
$ws.Range("A168:V169").PasteSpecial(-4122)
$excel.CutCopyMode = 0
